$d = $word.ActiveDocument

# Locate the "Docente(s) Responsável(eis)" heading paragraph.
$target = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    $t = $t.TrimEnd([char]13, [char]7)
    if ($t -eq "Docente(s) Responsável(eis) ") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    Write-Host "Target paragraph 'Docente(s) Responsável(eis)' not found."
} else {
    # Insert a brand-new paragraph right after it.
    $target.Range.InsertParagraphAfter() | Out-Null
    $newPara = $target.Next()

    # Give it the "List Bullet" style and the teacher's name/ID as its text.
    $newPara.Style = "ListBullet"
    $newPara.Range.Text = "6270264 - Juan Fernando Zapata Zapata"

    Write-Host "Inserted paragraph: $($newPara.Range.Text)"
}
